# Scheduled runner update: refresh Universalis market-price columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the Hades_Profits leve-crafting sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1054.375
$ws.Range("J40").Value = 1074.0625
$ws.Range("L40").Value = 1074.0625
$ws.Range("N40").Value = -1424.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1418.85
$ws.Range("I2").Value = 1464.2667
$ws.Range("J2").Value = 1282.6
$ws.Range("K2").Value = 1464.2667
$ws.Range("L2").Value = 1282.6
$ws.Range("M2").Value = -1351.2667
$ws.Range("N2").Value = -1508.6

$ws.Range("H61").Value = 38540250
$ws.Range("I61").Value = 52685856
$ws.Range("J61").Value = 145029.14
$ws.Range("K61").Value = 52685856
$ws.Range("L61").Value = 145029.14
$ws.Range("M61").Value = -52685644
$ws.Range("N61").Value = -145453.14

$ws.Range("H106").Value = 41316.668
$ws.Range("J106").Value = 41316.668
$ws.Range("L106").Value = 41316.668
$ws.Range("N106").Value = -43840.668

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 39401
$ws.Range("J113").Value = 39401
$ws.Range("L113").Value = 39401
$ws.Range("N113").Value = -48079

$ws.Range("H116").Value = 1418.85
$ws.Range("I116").Value = 1464.2667
$ws.Range("J116").Value = 1282.6
$ws.Range("K116").Value = 1464.2667
$ws.Range("L116").Value = 1282.6
$ws.Range("M116").Value = 829.7333000000001
$ws.Range("N116").Value = -5870.6

$ws.Range("H136").Value = 38540250
$ws.Range("I136").Value = 52685856
$ws.Range("J136").Value = 145029.14
$ws.Range("K136").Value = 158057568
$ws.Range("L136").Value = 435087.42
$ws.Range("M136").Value = -158055018
$ws.Range("N136").Value = -440187.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1418.85
$ws.Range("I3").Value = 1464.2667
$ws.Range("J3").Value = 1282.6
$ws.Range("K3").Value = 1464.2667
$ws.Range("L3").Value = 1282.6
$ws.Range("M3").Value = -1350.2667
$ws.Range("N3").Value = -1510.6

$ws.Range("H22").Value = 809.55554
$ws.Range("I22").Value = 574.5
$ws.Range("J22").Value = 997.6
$ws.Range("K22").Value = 574.5
$ws.Range("L22").Value = 997.6
$ws.Range("M22").Value = -401.5
$ws.Range("N22").Value = -1343.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8999.714
$ws.Range("I4").Value = 3000
$ws.Range("J4").Value = 9999.666999999999
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 9999.666999999999
$ws.Range("M4").Value = -2888
$ws.Range("N4").Value = -10223.667

$ws.Range("H6").Value = 40001000
$ws.Range("I6").Value = 80000000
$ws.Range("K6").Value = 80000000
$ws.Range("M6").Value = -79999887

$ws.Range("H7").Value = 116.3125
$ws.Range("I7").Value = 38.9
$ws.Range("J7").Value = 245.33333
$ws.Range("K7").Value = 38.9
$ws.Range("L7").Value = 245.33333
$ws.Range("M7").Value = 74.09999999999999
$ws.Range("N7").Value = -471.33333

$ws.Range("H31").Value = 1464.1818
$ws.Range("I31").Value = 1300.5807
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1300.5807
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1005.5807
$ws.Range("N31").Value = -4590

$ws.Range("H34").Value = 1464.1818
$ws.Range("I34").Value = 1300.5807
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1300.5807
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1098.5807
$ws.Range("N34").Value = -4404

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12472.454
$ws.Range("I80").Value = 26975
$ws.Range("J80").Value = 4185.2856
$ws.Range("K80").Value = 26975
$ws.Range("L80").Value = 4185.2856
$ws.Range("M80").Value = -25977
$ws.Range("N80").Value = -6181.2856

$ws.Range("H83").Value = 12472.454
$ws.Range("I83").Value = 26975
$ws.Range("J83").Value = 4185.2856
$ws.Range("K83").Value = 134875
$ws.Range("L83").Value = 20926.428
$ws.Range("M83").Value = -129883
$ws.Range("N83").Value = -30910.428

$ws.Range("H100").Value = 33560
$ws.Range("J100").Value = 33560
$ws.Range("L100").Value = 33560
$ws.Range("N100").Value = -35724

$ws.Range("H101").Value = 43630.25
$ws.Range("J101").Value = 43630.25
$ws.Range("L101").Value = 43630.25
$ws.Range("N101").Value = -50120.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2600.5
$ws.Range("I61").Value = 2569.7693
$ws.Range("K61").Value = 2569.7693
$ws.Range("M61").Value = -2367.7693

$ws.Range("H76").Value = 37063.438
$ws.Range("J76").Value = 37063.438
$ws.Range("L76").Value = 37063.438
$ws.Range("N76").Value = -37739.438

$ws.Range("H79").Value = 37063.438
$ws.Range("J79").Value = 37063.438
$ws.Range("L79").Value = 37063.438
$ws.Range("N79").Value = -39403.438

$ws.Range("H82").Value = 71431240
$ws.Range("I82").Value = 333334750
$ws.Range("J82").Value = 3000.818
$ws.Range("K82").Value = 333334750
$ws.Range("L82").Value = 3000.818
$ws.Range("M82").Value = -333334389
$ws.Range("N82").Value = -3722.818

$ws.Range("H85").Value = 71431240
$ws.Range("I85").Value = 333334750
$ws.Range("J85").Value = 3000.818
$ws.Range("K85").Value = 333334750
$ws.Range("L85").Value = 3000.818
$ws.Range("M85").Value = -333333502
$ws.Range("N85").Value = -5496.818

$ws.Range("H103").Value = 31250
$ws.Range("J103").Value = 31250
$ws.Range("L103").Value = 31250
$ws.Range("N103").Value = -33594

$ws.Range("H113").Value = 2600.5
$ws.Range("I113").Value = 2569.7693
$ws.Range("K113").Value = 2569.7693
$ws.Range("M113").Value = -399.7692999999999

$ws.Range("H136").Value = 145259.58
$ws.Range("I136").Value = 113822.555
$ws.Range("J136").Value = 201846.2
$ws.Range("K136").Value = 341467.665
$ws.Range("L136").Value = 605538.6000000001
$ws.Range("M136").Value = -338917.665
$ws.Range("N136").Value = -610638.6000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 6000
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6298

$ws.Range("H117").Value = 49608
$ws.Range("J117").Value = 49608
$ws.Range("L117").Value = 49608
$ws.Range("N117").Value = -58786
